$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -19.70661573854327
$ws.Range("C2").Value = 1.965854925203151
$ws.Range("D2").Value = -19.70661573854327
$ws.Range("E2").Value = -19.70661573854327
$ws.Range("F2").Value = -19.70661573854327
$ws.Range("G2").Value = -19.70661573854327
$ws.Range("H2").Value = -19.70661573854327
$ws.Range("I2").Value = -19.70661573854327
$ws.Range("J2").Value = -19.70661573854327
$ws.Range("K2").Value = -19.70661573854327

$ws.Range("B3").Value = -19.70661573854327
$ws.Range("C3").Value = -19.70661573854327
$ws.Range("D3").Value = -19.70661573854327
$ws.Range("E3").Value = -19.70661573854327
$ws.Range("F3").Value = -19.70661573854327
$ws.Range("G3").Value = -19.70661573854327
$ws.Range("H3").Value = -19.70661573854327
$ws.Range("I3").Value = 1.272882725158964
$ws.Range("J3").Value = -19.70661573854327
$ws.Range("K3").Value = -19.70661573854327

$ws.Range("B4").Value = -19.70661573854327
$ws.Range("C4").Value = 1.989774964419976
$ws.Range("D4").Value = 1.675027013789923
$ws.Range("E4").Value = -19.70661573854327
$ws.Range("F4").Value = 3.475156638903425
$ws.Range("G4").Value = -19.70661573854327
$ws.Range("H4").Value = 1.390215853081835
$ws.Range("I4").Value = -19.70661573854327
$ws.Range("J4").Value = 1.231129472867302
$ws.Range("K4").Value = -19.70661573854327

$ws.Range("B5").Value = -19.70661573854327
$ws.Range("C5").Value = 1.676293729864971
$ws.Range("D5").Value = -19.70661573854327
$ws.Range("E5").Value = -19.70661573854327
$ws.Range("F5").Value = -19.70661573854327
$ws.Range("G5").Value = -19.70661573854327
$ws.Range("H5").Value = -19.70661573854327
$ws.Range("I5").Value = -19.70661573854327
$ws.Range("J5").Value = -19.70661573854327
$ws.Range("K5").Value = -19.70661573854327

$ws.Range("B6").Value = -19.70661573854327
$ws.Range("C6").Value = -19.70661573854327
$ws.Range("D6").Value = -19.70661573854327
$ws.Range("E6").Value = -19.70661573854327
$ws.Range("F6").Value = -19.70661573854327
$ws.Range("G6").Value = -19.70661573854327
$ws.Range("H6").Value = -19.70661573854327
$ws.Range("I6").Value = -19.70661573854327
$ws.Range("J6").Value = -19.70661573854327
$ws.Range("K6").Value = -19.70661573854327

$ws.Range("B7").Value = 2.450481043084002
$ws.Range("C7").Value = -19.70661573854327
$ws.Range("D7").Value = -19.70661573854327
$ws.Range("E7").Value = -19.70661573854327
$ws.Range("F7").Value = -19.70661573854327
$ws.Range("G7").Value = -19.70661573854327
$ws.Range("H7").Value = -19.70661573854327
$ws.Range("I7").Value = -19.70661573854327
$ws.Range("J7").Value = -19.70661573854327
$ws.Range("K7").Value = -19.70661573854327

$ws.Range("B8").Value = -19.70661573854327
$ws.Range("C8").Value = -19.70661573854327
$ws.Range("D8").Value = -19.70661573854327
$ws.Range("E8").Value = 1.787193300874039
$ws.Range("F8").Value = -19.70661573854327
$ws.Range("G8").Value = -19.70661573854327
$ws.Range("H8").Value = -19.70661573854327
$ws.Range("I8").Value = -19.70661573854327
$ws.Range("J8").Value = -19.70661573854327
$ws.Range("K8").Value = -19.70661573854327

$ws.Range("B9").Value = 3.861359479161995
$ws.Range("C9").Value = -19.70661573854327
$ws.Range("D9").Value = -19.70661573854327
$ws.Range("E9").Value = -19.70661573854327
$ws.Range("F9").Value = -19.70661573854327
$ws.Range("G9").Value = -19.70661573854327
$ws.Range("H9").Value = -19.70661573854327
$ws.Range("I9").Value = -19.70661573854327
$ws.Range("J9").Value = -19.70661573854327
$ws.Range("K9").Value = -19.70661573854327

$ws.Range("B10").Value = -19.70661573854327
$ws.Range("C10").Value = -19.70661573854327
$ws.Range("D10").Value = -19.70661573854327
$ws.Range("E10").Value = -19.70661573854327
$ws.Range("F10").Value = -19.70661573854327
$ws.Range("G10").Value = -19.70661573854327
$ws.Range("H10").Value = -19.70661573854327
$ws.Range("I10").Value = 1.722775167010648
$ws.Range("J10").Value = -19.70661573854327
$ws.Range("K10").Value = 1.958585270159671

$ws.Range("B11").Value = -19.70661573854327
$ws.Range("C11").Value = -19.70661573854327
$ws.Range("D11").Value = -19.70661573854327
$ws.Range("E11").Value = 2.936165998937112
$ws.Range("F11").Value = -19.70661573854327
$ws.Range("G11").Value = -19.70661573854327
$ws.Range("H11").Value = -19.70661573854327
$ws.Range("I11").Value = -19.70661573854327
$ws.Range("J11").Value = -19.70661573854327
$ws.Range("K11").Value = 2.011934278123509

$ws.Range("B12").Value = -19.70661573854327
$ws.Range("C12").Value = -19.70661573854327
$ws.Range("D12").Value = -19.70661573854327
$ws.Range("E12").Value = -19.70661573854327
$ws.Range("F12").Value = -19.70661573854327
$ws.Range("G12").Value = -19.70661573854327
$ws.Range("H12").Value = -19.70661573854327
$ws.Range("I12").Value = -19.70661573854327
$ws.Range("J12").Value = -19.70661573854327
$ws.Range("K12").Value = -19.70661573854327

$ws.Range("B13").Value = -19.70661573854327
$ws.Range("C13").Value = -19.70661573854327
$ws.Range("D13").Value = -19.70661573854327
$ws.Range("E13").Value = 2.501902333323789
$ws.Range("F13").Value = -19.70661573854327
$ws.Range("G13").Value = -19.70661573854327
$ws.Range("H13").Value = -19.70661573854327
$ws.Range("I13").Value = -19.70661573854327
$ws.Range("J13").Value = 1.75021039938113
$ws.Range("K13").Value = 1.91819265505959

$ws.Range("B14").Value = -19.70661573854327
$ws.Range("C14").Value = -19.70661573854327
$ws.Range("D14").Value = 1.524069100006189
$ws.Range("E14").Value = -19.70661573854327
$ws.Range("F14").Value = -19.70661573854327
$ws.Range("G14").Value = -19.70661573854327
$ws.Range("H14").Value = -19.70661573854327
$ws.Range("I14").Value = -19.70661573854327
$ws.Range("J14").Value = -19.70661573854327
$ws.Range("K14").Value = 2.127744587285671

$ws.Range("B15").Value = -19.70661573854327
$ws.Range("C15").Value = -19.70661573854327
$ws.Range("D15").Value = 1.731006816474995
$ws.Range("E15").Value = -19.70661573854327
$ws.Range("F15").Value = -19.70661573854327
$ws.Range("G15").Value = -19.70661573854327
$ws.Range("H15").Value = -19.70661573854327
$ws.Range("I15").Value = -19.70661573854327
$ws.Range("J15").Value = -19.70661573854327
$ws.Range("K15").Value = -19.70661573854327

$ws.Range("B16").Value = -19.70661573854327
$ws.Range("C16").Value = -19.70661573854327
$ws.Range("D16").Value = -19.70661573854327
$ws.Range("E16").Value = -19.70661573854327
$ws.Range("F16").Value = -19.70661573854327
$ws.Range("G16").Value = -19.70661573854327
$ws.Range("H16").Value = -19.70661573854327
$ws.Range("I16").Value = -19.70661573854327
$ws.Range("J16").Value = 1.856774672541239
$ws.Range("K16").Value = -19.70661573854327

$ws.Range("B17").Value = -19.70661573854327
$ws.Range("C17").Value = 2.175735132264677
$ws.Range("D17").Value = 1.873584392763839
$ws.Range("E17").Value = -19.70661573854327
$ws.Range("F17").Value = -19.70661573854327
$ws.Range("G17").Value = -19.70661573854327
$ws.Range("H17").Value = 2.04390193975739
$ws.Range("I17").Value = 1.946282175942906
$ws.Range("J17").Value = 2.442841452061024
$ws.Range("K17").Value = -19.70661573854327

$ws.Range("B18").Value = -19.70661573854327
$ws.Range("C18").Value = -19.70661573854327
$ws.Range("D18").Value = -19.70661573854327
$ws.Range("E18").Value = -19.70661573854327
$ws.Range("F18").Value = -19.70661573854327
$ws.Range("G18").Value = -19.70661573854327
$ws.Range("H18").Value = 2.127158031915799
$ws.Range("I18").Value = 2.057203260255406
$ws.Range("J18").Value = 2.386621593770588
$ws.Range("K18").Value = -19.70661573854327

$ws.Range("B19").Value = -19.70661573854327
$ws.Range("C19").Value = -19.70661573854327
$ws.Range("D19").Value = 2.059081916775032
$ws.Range("E19").Value = -19.70661573854327
$ws.Range("F19").Value = -19.70661573854327
$ws.Range("G19").Value = -19.70661573854327
$ws.Range("H19").Value = 1.688345630803952
$ws.Range("I19").Value = 1.924167291660627
$ws.Range("J19").Value = -19.70661573854327
$ws.Range("K19").Value = -19.70661573854327

$ws.Range("B20").Value = -19.70661573854327
$ws.Range("C20").Value = 0.9979203416279468
$ws.Range("D20").Value = 1.474823484236461
$ws.Range("E20").Value = -19.70661573854327
$ws.Range("F20").Value = 3.150466396930197
$ws.Range("G20").Value = -19.70661573854327
$ws.Range("H20").Value = 1.511408793898659
$ws.Range("I20").Value = 1.305720312122039
$ws.Range("J20").Value = -19.70661573854327
$ws.Range("K20").Value = 1.974573335784428

$ws.Range("B21").Value = -19.70661573854327
$ws.Range("C21").Value = 1.269305679727289
$ws.Range("D21").Value = -19.70661573854327
$ws.Range("E21").Value = 1.691752724577286
$ws.Range("F21").Value = -19.70661573854327
$ws.Range("G21").Value = 4.321926493059163
$ws.Range("H21").Value = 1.492057703549638
$ws.Range("I21").Value = -19.70661573854327
$ws.Range("J21").Value = -19.70661573854327
$ws.Range("K21").Value = -19.70661573854327
